$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column cells to stay plain text (matching the
# source inline-string cells) instead of being auto-parsed into
# numbers/dates when the new value is assigned below. NumberFormat
# is set per cell -- a multi-area "A1,B2" Range only applies to the
# first area in this engine, so each cell gets its own statement.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "29.038.14"
$ws.Range("E2").Value = "  -0.56%  "
$ws.Range("D3").Value = "1.829.87"
$ws.Range("E3").Value = "  -0.24%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "240.99"
$ws.Range("E5").Value = "  -0.71%  "
$ws.Range("D6").Value = "0.6229"
$ws.Range("E6").Value = "  -6.10%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "0.07586"
$ws.Range("E8").Value = "  +2.05%  "
$ws.Range("D9").Value = "44.53"
$ws.Range("E9").Value = "  +6.34%  "
$ws.Range("D10").Value = "0.2916"
$ws.Range("E10").Value = "  -0.90%  "
$ws.Range("D11").Value = "22.77"
$ws.Range("E11").Value = "  -0.59%  "
$ws.Range("D12").Value = "0.07623"
$ws.Range("E12").Value = "  -1.93%  "
$ws.Range("D13").Value = "1.832.37"
$ws.Range("E13").Value = "  +0.01%  "
$ws.Range("D14").Value = "4.957"
$ws.Range("E14").Value = "  -0.91%  "
$ws.Range("D15").Value = "0.6643"
$ws.Range("E15").Value = "  -0.52%  "
$ws.Range("D16").Value = "82.23"
$ws.Range("E16").Value = "  -1.00%  "
$ws.Range("D17").Value = "0.000009049"
$ws.Range("E17").Value = "  +7.59%  "
$ws.Range("D18").Value = "6.000"
$ws.Range("E18").Value = "  -1.90%  "
$ws.Range("D19").Value = "29.048.26"
$ws.Range("E19").Value = "  -0.41%  "
$ws.Range("D20").Value = "2.081.22"
$ws.Range("E20").Value = "  +0.83%  "
$ws.Range("D21").Value = "224.83"
$ws.Range("E21").Value = "  -2.05%  "
$ws.Range("D22").Value = "12.33"
$ws.Range("E22").Value = "  -1.27%  "
$ws.Range("D24").Value = "7.196"
$ws.Range("E24").Value = "  +0.64%  "
$ws.Range("D25").Value = "1.000"
$ws.Range("E25").Value = "  +0.05%  "
$ws.Range("D26").Value = "159.46"
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("D27").Value = "8.398"
$ws.Range("E27").Value = "  -2.68%  "
$ws.Range("D28").Value = "0.1360"
$ws.Range("E28").Value = "  -2.97%  "
$ws.Range("E29").Value = "  -1.11%  "
$ws.Range("D30").Value = "1.497"
$ws.Range("E30").Value = "  -1.08%  "
$ws.Range("D31").Value = "1.210"
$ws.Range("E31").Value = "  +1.46%  "
$ws.Range("D32").Value = "4.050"
$ws.Range("E32").Value = "  -1.80%  "
$ws.Range("D33").Value = "4.015"
$ws.Range("E33").Value = "  -0.93%  "
$ws.Range("D34").Value = "0.05218"
$ws.Range("E34").Value = "  -1.20%  "
$ws.Range("D35").Value = "1.841"
$ws.Range("E35").Value = "  -1.52%  "
$ws.Range("D36").Value = "1.155"
$ws.Range("E36").Value = "  +1.15%  "
$ws.Range("D37").Value = "0.7321"
$ws.Range("E37").Value = "  -1.77%  "
$ws.Range("D38").Value = "2.644"
$ws.Range("E38").Value = "  -0.29%  "
$ws.Range("D39").Value = "1.273.03"
$ws.Range("E39").Value = "  -3.27%  "
$ws.Range("D40").Value = "2.749"
$ws.Range("E40").Value = "  +0.45%  "
$ws.Range("D41").Value = "0.01780"
$ws.Range("E41").Value = "  -0.86%  "
$ws.Range("D42").Value = "6.369"
$ws.Range("E42").Value = "  +7.44%  "
$ws.Range("D43").Value = "0.8901"
$ws.Range("E43").Value = "  -4.32%  "
$ws.Range("D44").Value = "1.000"
$ws.Range("E44").Value = "  +0.10%  "
$ws.Range("D45").Value = "101.54"
$ws.Range("E45").Value = "  -1.21%  "
$ws.Range("D46").Value = "1.979.14"
$ws.Range("E46").Value = "  +0.46%  "
$ws.Range("D47").Value = "0.5117"
$ws.Range("E47").Value = "  -0.57%  "
$ws.Range("D48").Value = "63.46"
$ws.Range("E48").Value = "  +0.35%  "
$ws.Range("E49").Value = "  -1.10%  "
$ws.Range("D50").Value = "0.3965"
$ws.Range("E50").Value = "  -1.40%  "
$ws.Range("D51").Value = "8.853"
$ws.Range("E51").Value = "  +0.43%  "
